# Fix a typo in the "Accion" column: "TNGO4" -> "TGNO4"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A17").Value = "TGNO4"
